$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Market Cap (column C) values for rows whose Name/Symbol stay the same.
$ws.Range("C2").Value = 704448480175.6102
$ws.Range("C3").Value = 234301122105.401
$ws.Range("C4").Value = 34780481584.6476
$ws.Range("C5").Value = 31440446351.24937
$ws.Range("C6").Value = 22409126657.29711
$ws.Range("C7").Value = 12771770288.35733
$ws.Range("C8").Value = 10367684249.58487
$ws.Range("C9").Value = 8595727131.068867
$ws.Range("C10").Value = 7968598284.179463
$ws.Range("C11").Value = 7597837770.961227
$ws.Range("C12").Value = 6882296974.000452
$ws.Range("C13").Value = 6836765843.319241
$ws.Range("C14").Value = 6154115361.756376
$ws.Range("C15").Value = 5864055698.012118
$ws.Range("C16").Value = 4943884169.308591
$ws.Range("C17").Value = 4622630406.286113
$ws.Range("C18").Value = 4274405809.596792
$ws.Range("C19").Value = 3671686599.97737
$ws.Range("C20").Value = 3401380976.234488
$ws.Range("C21").Value = 3192916216.664364
$ws.Range("C22").Value = 3132420847.329648
$ws.Range("C23").Value = 2922229784.867324

# Rows 24 and 25 swap places (Kaspa now ranks above Ethereum Classic),
# and both get refreshed market cap figures.
$ws.Range("A24").Value = "Kaspa"
$ws.Range("B24").Value = "KAS-USD"
$ws.Range("C24").Value = 2678810291.723003

$ws.Range("A25").Value = "Ethereum Classic"
$ws.Range("B25").Value = "ETC-USD"
$ws.Range("C25").Value = 2632693755.895442

$ws.Range("C26").Value = 2193394215.989171
